# Apply updated cryptocurrency price/volume data (and two row re-orderings)
# to match the latest scrape, as produced by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + '37.739.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.09%  '

# Row 3
$ws.Range("D3").Value = "'" + '2.031.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = "'" + '227.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.09%  '

# Row 6
$ws.Range("D6").Value = "'" + '0.603'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '

# Row 7
$ws.Range("D7").Value = "'" + '59.99'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.66%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("E9").Value = '  -2.45%  '

# Row 10
$ws.Range("E10").Value = '  +3.09%  '

# Row 11
$ws.Range("E11").Value = '  -0.01%  '

# Row 12
$ws.Range("D12").Value = "'" + '14.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.32%  '

# Row 13
$ws.Range("D13").Value = "'" + '2.332.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.46%  '

# Row 14
$ws.Range("D14").Value = "'" + '21.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.19%  '

# Row 15
$ws.Range("D15").Value = "'" + '0.768'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.13%  '

# Row 16
$ws.Range("D16").Value = "'" + '5.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.44%  '

# Row 17
$ws.Range("D17").Value = "'" + '2.038.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.28%  '

# Row 18
$ws.Range("D18").Value = "'" + '37.714.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.17%  '

# Row 19
$ws.Range("D19").Value = "'" + '69.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.67%  '

# Row 20
$ws.Range("D20").Value = "'" + '5.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.85%  '

# Row 21
$ws.Range("D21").Value = "'" + '0.0₃0825'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.79%  '

# Row 22
$ws.Range("D22").Value = "'" + '223.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.64%  '

# Row 23
$ws.Range("E23").Value = '  +0.23%  '

# Row 24
$ws.Range("D24").Value = "'" + '2.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.31%  '

# Row 25
$ws.Range("D25").Value = "'" + '2.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.69%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = "'" + '9.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.40%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'" + '167.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.84%  '

# Row 28
$ws.Range("D28").Value = "'" + '0.130'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.53%  '

# Row 29
$ws.Range("D29").Value = "'" + '18.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.65%  '

# Row 30
$ws.Range("E30").Value = '  -2.97%  '

# Row 31
$ws.Range("E31").Value = '  +0.87%  '

# Row 32
$ws.Range("E32").Value = '  +10.33%  '

# Row 33
$ws.Range("D33").Value = "'" + '4.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.68%  '

# Row 34
$ws.Range("D34").Value = "'" + '0.0606'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.57%  '

# Row 35
$ws.Range("D35").Value = "'" + '4.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.17%  '

# Row 36
$ws.Range("D36").Value = "'" + '6.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.63%  '

# Row 37
$ws.Range("D37").Value = "'" + '2.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.54%  '

# Row 38
$ws.Range("E38").Value = '  +3.72%  '

# Row 39
$ws.Range("E39").Value = '  -0.05%  '

# Row 40
$ws.Range("D40").Value = "'" + '18.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.91%  '

# Row 41
$ws.Range("D41").Value = "'" + '1.535.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.25%  '

# Row 42
$ws.Range("E42").Value = '  -0.23%  '

# Row 43
$ws.Range("D43").Value = "'" + '96.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.57%  '

# Row 44
$ws.Range("E44").Value = '  -2.71%  '

# Row 45
$ws.Range("D45").Value = "'" + '0.0909'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.26%  '

# Row 46
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = "'" + '1.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.05%  '

# Row 47
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = "'" + '4.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.19%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = "'" + '1.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.06%  '

# Row 49
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = "'" + '2.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.36%  '

# Row 50
$ws.Range("D50").Value = "'" + '7.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '

# Row 51
$ws.Range("D51").Value = "'" + '2.221.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.50%  '
